# SQE_holdings.xlsx update:
#  - bump the "as of" date in the confidential disclaimer footnote (A33)
#  - refresh the Weight (D) / Percent Change (E) values for every holding row (2-30)
#
# The worksheet ships with sheetProtection enabled, so it must be
# unprotected before writing and re-protected afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# --- Footnote text: date 2021-06-10 -> 2021-06-14 -------------------------
$disclaimer = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-06-14 for illustrative purposes only and are subject to change."
$ws.Range("A33").Value = $disclaimer

# --- Weight (D) / Percent Change (E) refresh, rows 2-30 --------------------
$ws.Range("D2").Value = 0.004824501144123577
$ws.Range("E2").Value = 0.00518606132290067

$ws.Range("D3").Value = 0.01397565649578707
$ws.Range("E3").Value = 0.02898791708236326

$ws.Range("D4").Value = 0.3074711562626142
$ws.Range("E4").Value = 0.007698954818533554

$ws.Range("D5").Value = 0.3283851687418638
$ws.Range("E5").Value = 0.01106718895193359

$ws.Range("D6").Value = 0.01775656655874776
$ws.Range("E6").Value = 0.02457793482528459

$ws.Range("D7").Value = 0.001514119825800824
$ws.Range("E7").Value = -0.00784447476125516

$ws.Range("D8").Value = 0.003242550609209882
$ws.Range("E8").Value = -0.01194457716196851

$ws.Range("D9").Value = 0.003394065874179173
$ws.Range("E9").Value = -0.01095490231878771

$ws.Range("D10").Value = 0.002900169488982751
$ws.Range("E10").Value = -0.01086182336182318

$ws.Range("D11").Value = 0.003231086264007707
$ws.Range("E11").Value = 0.007032348804500765

$ws.Range("D12").Value = 0.01710666212465146
$ws.Range("E12").Value = 0.01663346012195865

$ws.Range("D13").Value = 0.03289949479677524
$ws.Range("E13").Value = -0.004665830035074259

$ws.Range("D14").Value = 0.002987443107864177
$ws.Range("E14").Value = 0.005877268798617052

$ws.Range("D15").Value = 0.01533371663149342
$ws.Range("E15").Value = 0.002485451018428675

$ws.Range("D16").Value = 0.01158859391639358
$ws.Range("E16").Value = -0.01696924324661553

$ws.Range("D17").Value = 0.03680597042441652
$ws.Range("E17").Value = 0.00279069767441853

$ws.Range("D18").Value = 0.05926397716054559
$ws.Range("E18").Value = 0.007755244484082269

$ws.Range("D19").Value = 0.007494893137714079
$ws.Range("E19").Value = -0.004395937547369888

$ws.Range("D20").Value = 0.02209241289892726
$ws.Range("E20").Value = 0.01085538772247241

$ws.Range("D21").Value = 0.00428105503270424
$ws.Range("E21").Value = -0.003377563329312383

$ws.Range("D22").Value = 0.00560720091014507
$ws.Range("E22").Value = 0.01282004052311669

$ws.Range("D23").Value = 0.001347706076188165
$ws.Range("E23").Value = -0.01295143212951411

$ws.Range("D24").Value = 0.009750064106087954
$ws.Range("E24").Value = 0.0003707548568885333

$ws.Range("D25").Value = 0.006204947737690928
$ws.Range("E25").Value = 0.0247596854063501

$ws.Range("D26").Value = 0.03287602387382755
$ws.Range("E26").Value = 0.003191836939857806

$ws.Range("D27").Value = 0.003108619171003387
$ws.Range("E27").Value = -0.002616431187859813

$ws.Range("D28").Value = 0.02669395318534637
$ws.Range("E28").Value = -0.003745318352059823

$ws.Range("D29").Value = 0.01786222444290835
$ws.Range("E29").Value = 0.00451009132934943

$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 0.007805359999675909

# Restore sheet protection to match the original state.
$ws.Protect()
